$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (price and volume%) to be treated as text so that
# values like "0.4500", "1.001", "30.106.92" keep their exact original
# formatting instead of being auto-converted into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.106.92'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '1.907.85'
$ws.Range('E3').Value = '  -1.36%  '
$ws.Range('D5').Value = '0.7424'
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('D6').Value = '244.05'
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.3085'
$ws.Range('E8').Value = '  -3.20%  '
$ws.Range('D9').Value = '26.51'
$ws.Range('E9').Value = '  -6.56%  '
$ws.Range('D10').Value = '0.06980'
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('D11').Value = '0.08077'
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').Value = '0.7683'
$ws.Range('E12').Value = '  -1.99%  '
$ws.Range('D13').Value = '1.932.53'
$ws.Range('E13').Value = '  -0.14%  '
$ws.Range('D14').Value = '5.319'
$ws.Range('E14').Value = '  -1.58%  '
$ws.Range('D15').Value = '92.34'
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('D16').Value = '14.26'
$ws.Range('E16').Value = '  -1.75%  '
$ws.Range('D17').Value = '30.111.78'
$ws.Range('E17').Value = '  -0.68%  '
$ws.Range('D18').Value = '6.083'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('D19').Value = '0.000007822'
$ws.Range('E19').Value = '  -2.13%  '
$ws.Range('D20').Value = '240.05'
$ws.Range('E20').Value = '  -4.96%  '
$ws.Range('D21').Value = '2.181.45'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').Value = '7.131'
$ws.Range('E24').Value = '  +6.22%  '
$ws.Range('D25').Value = '9.389'
$ws.Range('E25').Value = '  -2.08%  '
$ws.Range('D26').Value = '167.22'
$ws.Range('E26').Value = '  +1.35%  '
$ws.Range('D27').Value = '18.98'
$ws.Range('E27').Value = '  -0.50%  '
$ws.Range('D28').Value = '0.1272'
$ws.Range('E28').Value = '  -2.70%  '
$ws.Range('D29').Value = '2.048'
$ws.Range('E29').Value = '  -7.36%  '
$ws.Range('D30').Value = '1.546'
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').Value = '1.353'
$ws.Range('E31').Value = '  -1.07%  '
$ws.Range('D32').Value = '4.336'
$ws.Range('E32').Value = '  -2.41%  '
$ws.Range('D33').Value = '4.084'
$ws.Range('E33').Value = '  -1.58%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.05232'
$ws.Range('E34').Value = '  -0.89%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.306'
$ws.Range('E35').Value = '  -2.18%  '
$ws.Range('D36').Value = '0.7495'
$ws.Range('E36').Value = '  -1.05%  '
$ws.Range('D37').Value = '2.723'
$ws.Range('E37').Value = '  -2.11%  '
$ws.Range('D38').Value = '0.01962'
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('D39').Value = '2.799'
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('E40').Value = '  -2.91%  '
$ws.Range('D41').Value = '0.4500'
$ws.Range('E41').Value = '  -0.37%  '
$ws.Range('D42').Value = '74.44'
$ws.Range('E42').Value = '  -5.94%  '
$ws.Range('D43').Value = '1.977'
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('D45').Value = '0.8408'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('D46').Value = '7.739'
$ws.Range('E46').Value = '  +0.29%  '
$ws.Range('D47').Value = '102.07'
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('D48').Value = '9.918'
$ws.Range('E48').Value = '  -1.84%  '
$ws.Range('D49').Value = '2.078.19'
$ws.Range('E49').Value = '  -0.41%  '
$ws.Range('D50').Value = '36.79'
$ws.Range('E50').Value = '  -2.41%  '
$ws.Range('D51').Value = '0.1182'
$ws.Range('E51').Value = '  -5.55%  '
